$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview" (sheet1) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").EntireRow.Insert()
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value2 = "0be54287-ee42-41f8-adc9-05dd855ab884.md"
$wsOverview.Range("B2").Value2 = "e2e\0be54287-ee42-41f8-adc9-05dd855ab884.md"
$wsOverview.Range("C2").Value2 = ".md"
$wsOverview.Range("D2").Value2 = ""
$wsOverview.Range("E2").Value2 = "Ready for handoff"
$wsOverview.Range("F2").Value2 = "Ready for handoff"
$wsOverview.Range("G2").Value2 = "2016-08-13 20:51:52"

# ---- Sheet "zh-cn" (sheet2) ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").EntireRow.Insert()
$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))
$wsZh.Hyperlinks.Delete()

$wsZh.Range("A2").Value2 = "0be54287-ee42-41f8-adc9-05dd855ab884.md"
$wsZh.Range("B2").Value2 = ".md"
$wsZh.Range("C2").Value2 = "Ready for handoff"
$wsZh.Range("D2").Value2 = "e2e"
$wsZh.Range("E2").Value2 = "ht"
$wsZh.Range("F2").Value2 = "False"
$wsZh.Range("G2").Value2 = "0be54287-ee42-41f8-adc9-05dd855ab884.f27a340706b68e066895494600b57b4006efe7a5.zh-cn.xlf"
$wsZh.Range("H2").Value2 = "2016-08-13 20:51:44"
$wsZh.Range("I2").Value2 = ""
$wsZh.Range("J2").Value2 = ""
$wsZh.Range("K2").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("L2").Value2 = ""
$wsZh.Range("M2").Value2 = "True"
$wsZh.Range("N2").Value2 = ""
$wsZh.Range("O2").Value2 = "False"
$wsZh.Range("P2").Value2 = ""

# ---- Sheet "de-de" (sheet3) ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").EntireRow.Insert()
$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
$wsDe.Hyperlinks.Delete()

$wsDe.Range("A2").Value2 = "0be54287-ee42-41f8-adc9-05dd855ab884.md"
$wsDe.Range("B2").Value2 = ".md"
$wsDe.Range("C2").Value2 = "Ready for handoff"
$wsDe.Range("D2").Value2 = "e2e"
$wsDe.Range("E2").Value2 = "ht"
$wsDe.Range("F2").Value2 = "False"
$wsDe.Range("G2").Value2 = "0be54287-ee42-41f8-adc9-05dd855ab884.f27a340706b68e066895494600b57b4006efe7a5.de-de.xlf"
$wsDe.Range("H2").Value2 = "2016-08-13 20:51:52"
$wsDe.Range("I2").Value2 = ""
$wsDe.Range("J2").Value2 = ""
$wsDe.Range("K2").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("L2").Value2 = ""
$wsDe.Range("M2").Value2 = "True"
$wsDe.Range("N2").Value2 = ""
$wsDe.Range("O2").Value2 = "False"
$wsDe.Range("P2").Value2 = ""

# ---- Hyperlinks: re-add in order so rId2 -> new file, rId3 -> old file ----
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/494604fe103f9887dcc9f35221e5cb7122d53a18/e2e/0be54287-ee42-41f8-adc9-05dd855ab884.md", "", "", "e2e\0be54287-ee42-41f8-adc9-05dd855ab884.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/494604fe103f9887dcc9f35221e5cb7122d53a18/e2e/79fd4040-5917-47f7-a65a-97dba0fe7202.md", "", "", "e2e\79fd4040-5917-47f7-a65a-97dba0fe7202.md")

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/494604fe103f9887dcc9f35221e5cb7122d53a18/e2e/0be54287-ee42-41f8-adc9-05dd855ab884.md", "", "", "0be54287-ee42-41f8-adc9-05dd855ab884.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/494604fe103f9887dcc9f35221e5cb7122d53a18/e2e/79fd4040-5917-47f7-a65a-97dba0fe7202.md", "", "", "79fd4040-5917-47f7-a65a-97dba0fe7202.md")

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/494604fe103f9887dcc9f35221e5cb7122d53a18/e2e/0be54287-ee42-41f8-adc9-05dd855ab884.md", "", "", "0be54287-ee42-41f8-adc9-05dd855ab884.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/494604fe103f9887dcc9f35221e5cb7122d53a18/e2e/79fd4040-5917-47f7-a65a-97dba0fe7202.md", "", "", "79fd4040-5917-47f7-a65a-97dba0fe7202.md")
